$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The percentage column holds plain text like "0.56%" / "0%" (not real
# numeric percentages). Simply assigning a string such as "0%" would make
# Excel auto-convert it into a numeric, percent-formatted cell, so we
# briefly force a text number format while writing the value, then put
# the cell style back to "Normal" (its original, unstyled state).
$ws.Range("B13").Value = "NTD 0"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "0%"
$ws.Range("C13").Style = "Normal"

$ws.Range("B14").Value = "NTD 0"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "0%"
$ws.Range("C14").Style = "Normal"

$ws.Range("B16").Value = "NTD 0"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "0%"
$ws.Range("C16").Style = "Normal"
